$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-typed cells for columns D and E (and B/C swap rows) so that
# numeric-looking strings like "1.00" or "83.60" are preserved exactly,
# matching the original inlineStr (text) cell type, instead of Excel
# auto-converting them to numbers.
$textRange = $ws.Range("D2:E51")
$textRange.NumberFormat = "@"

$ws.Range("D2").Value = "66.894.31"
$ws.Range("E2").Value = "  +3.70%  "

$ws.Range("D3").Value = "3.799.86"
$ws.Range("E3").Value = "  +4.22%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "421.57"
$ws.Range("E5").Value = "  +3.97%  "

$ws.Range("D6").Value = "128.96"
$ws.Range("E6").Value = "  -1.22%  "

$ws.Range("D7").Value = "3.798.36"
$ws.Range("E7").Value = "  +4.24%  "

$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  -2.47%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").Value = "0.718"
$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("D12").Value = "0.0000346"
$ws.Range("E12").Value = "  +10.89%  "

$ws.Range("D13").Value = "40.29"
$ws.Range("E13").Value = "  -3.44%  "

$ws.Range("D14").Value = "4.416.43"
$ws.Range("E14").Value = "  +4.76%  "

$ws.Range("D15").Value = "10.13"
$ws.Range("E15").Value = "  +3.13%  "

$ws.Range("D16").Value = "15.59"
$ws.Range("E16").Value = "  +18.48%  "

$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").Value = "3.793.24"
$ws.Range("E18").Value = "  +5.08%  "

$ws.Range("D19").Value = "19.55"
$ws.Range("E19").Value = "  -1.13%  "

$ws.Range("D20").Value = "66.986.01"
$ws.Range("E20").Value = "  +3.85%  "

$ws.Range("D21").Value = "1.08"
$ws.Range("E21").Value = "  +0.71%  "

$ws.Range("D22").Value = "404.81"
$ws.Range("E22").Value = "  -3.35%  "

$ws.Range("D23").Value = "15.05"
$ws.Range("E23").Value = "  -1.02%  "

$ws.Range("D24").Value = "83.60"
$ws.Range("E24").Value = "  -2.26%  "

$ws.Range("D25").Value = "3.02"
$ws.Range("E25").Value = "  +1.03%  "

$ws.Range("D26").Value = "37.03"
$ws.Range("E26").Value = "  +3.74%  "

$ws.Range("D27").Value = "5.47"
$ws.Range("E27").Value = "  +9.22%  "

$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("D29").Value = "9.70"
$ws.Range("E29").Value = "  +3.07%  "

$ws.Range("D30").Value = "9.09"
$ws.Range("E30").Value = "  +31.92%  "

$ws.Range("D31").Value = "732.58"
$ws.Range("E31").Value = "  +10.04%  "

$ws.Range("D32").Value = "12.48"
$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("D33").Value = "2.77"
$ws.Range("E33").Value = "  +2.68%  "

$ws.Range("D34").Value = "0.121"
$ws.Range("E34").Value = "  +2.63%  "

$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("E36").Value = "  -4.10%  "

$ws.Range("D37").Value = "38.47"
$ws.Range("E37").Value = "  -4.04%  "

$ws.Range("D38").Value = "55.06"
$ws.Range("E38").Value = "  -1.34%  "

$ws.Range("D39").Value = "5.36"
$ws.Range("E39").Value = "  +24.32%  "

$ws.Range("D40").Value = "0.0₃0740"
$ws.Range("E40").Value = "  +18.08%  "

$ws.Range("D41").Value = "0.0450"
$ws.Range("E41").Value = "  -2.26%  "

$ws.Range("D42").Value = "2.92"
$ws.Range("E42").Value = "  -0.57%  "

$ws.Range("E43").Value = "  +0.87%  "

$ws.Range("E44").Value = "  -3.94%  "

$ws.Range("E45").Value = "  +0.32%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.11"
$ws.Range("E46").Value = "  +0.99%  "

$ws.Range("D47").Value = "143.39"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("B48").Value = "TheGraph"
$ws.Range("C48").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D48").Value = "0.311"
$ws.Range("E48").Value = "  +7.94%  "

$ws.Range("E49").Value = "  -1.05%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "2.78"
$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("B51").Value = "WEMIXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "2.55"
$ws.Range("E51").Value = "  +1.85%  "

# Restore default (General) formatting on the cells we text-formatted,
# without disturbing the values that were just entered as text.
$textRange.ClearFormats()

